$d = $word.ActiveDocument

# --- IBM experience bullet 1: "Converted internal build tools to production" ---
# becomes: "Converted internal scripts to production quality to get integrated into WebSphere build tools. "
$d.Content.Find.Execute(
    "Converted internal build tools to production",
    $false, $true, $false, $false, $false, $true, 1, $false,
    "Converted internal scripts to production quality to get integrated into WebSphere build tools. ",
    2
) | Out-Null

# --- IBM experience bullet 2: Github Enterprise Perl library sentence ---
# "adding git to supported source control management"
# becomes: "following existing build structure, adding git to supported version control systems "
$d.Content.Find.Execute(
    "adding git to supported source control management",
    $false, $true, $false, $false, $false, $true, 1, $false,
    "following existing build structure, adding git to supported version control systems ",
    2
) | Out-Null

# The _GoBack bookmark moves from the end of the Github bullet to mid-way through
# the "Converted..." bullet (right after "erted internal scripts"). Recreate it there.
$r = $d.Content
$r.Find.Execute("Converted internal scripts", $false, $true, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
if ($r.Find.Found) {
    $gobackRange = $d.Range($r.End, $r.End)
    $d.Bookmarks.Add("_GoBack", $gobackRange) | Out-Null
}
